# Fixed some bugs in PiggyBank
# The underlying data rows (A2:F23) got reshuffled into a different row
# order. Snapshot the current A:F values for rows 2-23, then write them
# back out in the new row order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 23

# Snapshot current values (A:F) for every data row before moving anything,
# so writes below don't clobber data we still need to read.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le 6; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# For each new row (2..23), which old row's data should land there.
$newRowSource = @{
    2  = 8
    3  = 6
    4  = 13
    5  = 14
    6  = 10
    7  = 15
    8  = 3
    9  = 4
    10 = 11
    11 = 5
    12 = 9
    13 = 12
    14 = 2
    15 = 7
    16 = 19
    17 = 16
    18 = 17
    19 = 20
    20 = 21
    21 = 18
    22 = 23
    23 = 22
}

foreach ($newRow in $newRowSource.Keys) {
    $oldRow = $newRowSource[$newRow]
    $vals = $snapshot[$oldRow]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($newRow, $c).Value = $vals[$c - 1]
    }
}
